# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> @{ Col = Value } for cells that changed.
# Rows 50 and 51 also swap Coin name / Link (BitcoinSV <-> Algorand entries).
$updates = [ordered]@{
    2  = @{ D = "41.126.51";  E = "  -3.69%  " }
    3  = @{ D = "2.452.32";   E = "  -3.01%  " }
    4  = @{ E = "  +0.06%  " }
    5  = @{ D = "308.87";     E = "  -0.09%  " }
    6  = @{ D = "92.52";      E = "  -8.25%  " }
    7  = @{ E = "  -3.03%  " }
    8  = @{ E = "  +0.06%  " }
    9  = @{ D = "0.491";      E = "  -6.36%  " }
    10 = @{ D = "32.96";      E = "  -7.96%  " }
    11 = @{ D = "0.0770";     E = "  -4.27%  " }
    12 = @{ E = "  -1.18%  " }
    13 = @{ D = "6.91";       E = "  -5.97%  " }
    14 = @{ D = "2.836.41";   E = "  -2.75%  " }
    15 = @{ D = "2.453.89";   E = "  -4.34%  " }
    16 = @{ D = "14.52";      E = "  -5.69%  " }
    17 = @{ D = "0.772";      E = "  -4.72%  " }
    18 = @{ D = "41.122.77";  E = "  -3.67%  " }
    19 = @{ D = "6.24";       E = "  -7.22%  " }
    20 = @{ D = "0.0₃0910";   E = "  -4.43%  " }
    21 = @{ D = "11.10";      E = "  -10.06%  " }
    22 = @{ D = "67.34";      E = "  -3.35%  " }
    23 = @{ D = "234.13";     E = "  -4.10%  " }
    24 = @{ D = "2.73";       E = "  -5.01%  " }
    25 = @{ E = "  +0.28%  " }
    26 = @{ D = "1.88";       E = "  -7.65%  " }
    27 = @{ D = "23.65";      E = "  -7.26%  " }
    28 = @{ E = "  -5.70%  " }
    29 = @{ D = "9.51";       E = "  -6.11%  " }
    30 = @{ D = "35.16";      E = "  -9.24%  " }
    31 = @{ D = "150.70";     E = "  -4.28%  " }
    32 = @{ D = "5.41";       E = "  -5.77%  " }
    33 = @{ D = "2.66";       E = "  -5.13%  " }
    34 = @{ E = "  -3.14%  " }
    35 = @{ D = "0.0728";     E = "  -7.27%  " }
    36 = @{ D = "2.95";       E = "  -6.68%  " }
    37 = @{ D = "1.84";       E = "  -7.53%  " }
    38 = @{ D = "16.56";      E = "  -7.10%  " }
    39 = @{ E = "  -4.08%  " }
    40 = @{ D = "0.101";      E = "  -9.07%  " }
    41 = @{ D = "4.09";       E = "  -2.40%  " }
    42 = @{ E = "  +0.20%  " }
    43 = @{ D = "19.89";      E = "  -8.83%  " }
    44 = @{ D = "1.967.03";   E = "  -2.06%  " }
    45 = @{ D = "0.0280";     E = "  -6.81%  " }
    46 = @{ D = "2.97";       E = "  -9.55%  " }
    47 = @{ D = "8.48";       E = "  -4.66%  " }
    48 = @{ D = "70.00";      E = "  -3.09%  " }
    49 = @{ D = "95.47";      E = "  -5.83%  " }
    50 = @{ B = "BitcoinSV"; C = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"; D = "73.21";  E = "  -7.50%  " }
    51 = @{ B = "Algorand";  C = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D = "0.174";  E = "  -8.25%  " }
}

# Cells whose new text would otherwise be silently re-interpreted by Excel as a
# number and lose formatting (e.g. trailing/leading zeros: "11.10" -> 11.1,
# "70.00" -> 70, "0.0280" -> 0.028). Force these to text first so the literal
# string is preserved, then drop back to the default "Normal" style so no
# stray number-format style is left behind on the cell.
$textForceCells = @("D11", "D21", "D31", "D45", "D48")

foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

foreach ($rowNum in $updates.Keys) {
    $rowChanges = $updates[$rowNum]
    foreach ($col in $rowChanges.Keys) {
        $ws.Range("$col$rowNum").Value = $rowChanges[$col]
    }
}

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
